# Add a new worksheet "population_5_countries" right after "population_benelux"
# and before "births", containing population data (Male/Female, 2013-2017) for
# Belgium, France, Germany, Luxembourg and Netherlands.

$wb = $excel.ActiveWorkbook

$afterSheet = $wb.Worksheets.Item("population_benelux")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$newSheet.Name = "population_5_countries"

$data = @(
    @("country", "gender\time", 2013, 2014, 2015, 2016, 2017),
    @("Belgium", "Male", 5472856, 5493792, 5524068, 5569264, 5589272),
    @("Belgium", "Female", 5665118, 5687048, 5713206, 5741853, 5762455),
    @("France", "Male", 31772665, 32045129, 32174258, 32247386, 32318973),
    @("France", "Female", 33827685, 34120851, 34283895, 34391005, 34485148),
    @("Germany", "Male", 39380976, 39556923, 39835457, 40514123, 40697118),
    @("Germany", "Female", 41142770, 41210540, 41362080, 41661561, 41824535),
    @("Luxembourg", "Male", 268412, 275117, 281972, 289193, 296641),
    @("Luxembourg", "Female", 268627, 274563, 280986, 287056, 294026),
    @("Netherlands", "Male", 8307339, 8334385, 8372858, 8417135, 8475102),
    @("Netherlands", "Female", 8472236, 8494904, 8527868, 8561985, 8606405)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $newSheet.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

# Restore the originally-active sheet ("immigration") as the selected tab,
# since adding a new sheet makes it active by default.
$immigration = $wb.Worksheets.Item("immigration")
$immigration.Activate()
